# Nexial "localdb" command-type addition
# -----------------------------------------------------------------
# The hidden "#system" sheet keeps one lookup column per command-type:
#   - column A ("target") lists every command-type name, alphabetically
#   - every other column is headed by a command-type name and lists
#     that command-type's individual command signatures underneath
#
# This script inserts a brand-new column (N) for the new "localdb"
# command-type (pushing every column from N onward one slot to the
# right), fills in its header + six command signatures, inserts
# "localdb" into the alphabetical target list in column A, and fixes
# up every defined name whose range shifted as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------
# 1. Make room for the new "localdb" column by inserting a column
#    at N - this pushes the existing N..AC columns (macro..xml) one
#    column to the right (O..AD).
# ---------------------------------------------------------------
$ws.Range("N1").EntireColumn.Insert()

# ---------------------------------------------------------------
# 2. Populate the new "localdb" column (header + 6 commands)
# ---------------------------------------------------------------
$ws.Range("N1").Value2 = "localdb"
$ws.Range("N2").Value2 = "cloneTable(var,source,target)"
$ws.Range("N3").Value2 = "dropTables(var,tables)"
$ws.Range("N4").Value2 = "exportCSV(sql,output)"
$ws.Range("N5").Value2 = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value2 = "purge(var)"
$ws.Range("N7").Value2 = "runSQLs(var,sqls)"

# ---------------------------------------------------------------
# 3. Insert "localdb" into the alphabetical "target" list in column
#    A (it sits between "json" and "macro"), pushing A14:A29 down to
#    A15:A30.
# ---------------------------------------------------------------
$ws.Range("A14").Insert()
$ws.Range("A14").Value2 = "localdb"

# ---------------------------------------------------------------
# 4. Fix up the defined names whose single-column ranges shifted one
#    column to the right because of the column insert in step 1.
# ---------------------------------------------------------------
$wb.Names.Item("macro").RefersTo      = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo       = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo     = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo        = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo      = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo      = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo        = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo      = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo        = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo       = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$21"

# ---------------------------------------------------------------
# 5. "target" (column A) now has one more row of data (29 -> 30)
# ---------------------------------------------------------------
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"

# ---------------------------------------------------------------
# 6. Register the new "localdb" defined name
# ---------------------------------------------------------------
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
